# Applies the "bug fixes" commit to Sablefish_Inputs.xlsx:
#  - Controls!B2: n_sims 1 -> 20
#  - Rename sheet "Recruitment" -> "Recruitment_Mortality"
#  - Recruitment_Mortality!B5: r0 1.2 -> 0.75
#  - Recruitment_Mortality!A6:C6: new "M" / 0.1 / "Natural Mortality" row
#  - Recruitment_Mortality!A7:C7: new "mu_rec" / 16.5 / "Mean recruitment" row
#  - Update the active selections left behind on a couple of sheets

$wb = $excel.ActiveWorkbook

# --- Controls sheet ---------------------------------------------------
$wsControls = $wb.Worksheets.Item("Controls")
$wsControls.Range("B2").Value = 20
$wsControls.Range("C3").Select()

# --- Recruitment sheet: rename + new Natural-Mortality / mean-recruit rows
$wsRecruit = $wb.Worksheets.Item("Recruitment")
$wsRecruit.Name = "Recruitment_Mortality"

$wsRecruit.Range("B5").Value = 0.75

$wsRecruit.Range("A6").Value = "M"
$wsRecruit.Range("B6").Value = 0.1
$wsRecruit.Range("C6").Value = "Natural Mortality"

$wsRecruit.Range("C7").Value = "Mean recruitment"
$wsRecruit.Range("A7").Value = "mu_rec"
$wsRecruit.Range("B7").Value = 16.5

$wsRecruit.Range("B6").Select()

$wb.Save()
